$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3310
$ws.Range("I40").Value = 2041
$ws.Range("K40").Value = 2041
$ws.Range("M40").Value = -1866
$ws.Range("H80").Value = 1932.4147
$ws.Range("J80").Value = 2217.28
$ws.Range("L80").Value = 6651.84
$ws.Range("N80").Value = -8647.84
$ws.Range("H83").Value = 1932.4147
$ws.Range("J83").Value = 2217.28
$ws.Range("L83").Value = 19955.52
$ws.Range("N83").Value = -29939.52
$ws.Range("H97").Value = 1184.5
$ws.Range("J97").Value = 1184.5
$ws.Range("L97").Value = 3553.5
$ws.Range("N97").Value = -4545.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11068.394
$ws.Range("I32").Value = 10115.862
$ws.Range("K32").Value = 10115.862
$ws.Range("M32").Value = -9828.861999999999
$ws.Range("H45").Value = 2811.4443
$ws.Range("I45").Value = 2159.75
$ws.Range("J45").Value = 3332.8
$ws.Range("K45").Value = 2159.75
$ws.Range("L45").Value = 3332.8
$ws.Range("M45").Value = -1782.75
$ws.Range("N45").Value = -4086.8
$ws.Range("H61").Value = 3596.2727
$ws.Range("I61").Value = 1495
$ws.Range("K61").Value = 1495
$ws.Range("M61").Value = -1283
$ws.Range("H97").Value = 1626.6154
$ws.Range("I97").Value = 1512.1666
$ws.Range("K97").Value = 1512.1666
$ws.Range("M97").Value = -1016.1666
$ws.Range("H110").Value = 1417.5
$ws.Range("I110").Value = 1417.5
$ws.Range("K110").Value = 1417.5
$ws.Range("M110").Value = 627.5
$ws.Range("H122").Value = 3145.5
$ws.Range("I122").Value = 2043.7667
$ws.Range("K122").Value = 6131.300099999999
$ws.Range("M122").Value = -3681.300099999999
$ws.Range("H136").Value = 3596.2727
$ws.Range("I136").Value = 1495
$ws.Range("K136").Value = 4485
$ws.Range("M136").Value = -1935

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3365.4814
$ws.Range("I20").Value = 2210.25
$ws.Range("K20").Value = 2210.25
$ws.Range("M20").Value = -1963.25
$ws.Range("H94").Value = 3955.2
$ws.Range("I94").Value = 3803.625
$ws.Range("K94").Value = 3803.625
$ws.Range("M94").Value = -3352.625
$ws.Range("H99").Value = 47007.11
$ws.Range("I99").Value = 52406.625
$ws.Range("K99").Value = 52406.625
$ws.Range("M99").Value = -50908.625
$ws.Range("H105").Value = 3498.75
$ws.Range("I105").Value = 3998.3333
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 3998.3333
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -2251.3333
$ws.Range("N105").Value = -5494
$ws.Range("H120").Value = 11999
$ws.Range("J120").Value = 11999
$ws.Range("L120").Value = 11999
$ws.Range("N120").Value = -21675
$ws.Range("H131").Value = 14216.454
$ws.Range("I131").Value = 10709
$ws.Range("J131").Value = 30000
$ws.Range("K131").Value = 10709
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080
$ws.Range("M131").Value = -5669
$ws.Range("H134").Value = 3344.0476
$ws.Range("I134").Value = 2045.6562
$ws.Range("K134").Value = 6136.9686
$ws.Range("M134").Value = -3601.9686

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 36261
$ws.Range("I86").Value = 51424.855
$ws.Range("K86").Value = 51424.855
$ws.Range("M86").Value = -50301.855
$ws.Range("H89").Value = 36261
$ws.Range("I89").Value = 51424.855
$ws.Range("K89").Value = 257124.275
$ws.Range("M89").Value = -251508.275
$ws.Range("H99").Value = 11224083
$ws.Range("I99").Value = 2034221.1
$ws.Range("K99").Value = 2034221.1
$ws.Range("M99").Value = -2032723.1
$ws.Range("H107").Value = 4295.593
$ws.Range("I107").Value = 605.2174
$ws.Range("J107").Value = 25515.25
$ws.Range("K107").Value = 605.2174
$ws.Range("L107").Value = 25515.25
$ws.Range("M107").Value = 1314.7826
$ws.Range("N107").Value = -29355.25
$ws.Range("H126").Value = 11224083
$ws.Range("I126").Value = 2034221.1
$ws.Range("K126").Value = 6102663.300000001
$ws.Range("M126").Value = -6100193.300000001
$ws.Range("H134").Value = 3243.25
$ws.Range("I134").Value = 1901.6487
$ws.Range("K134").Value = 5704.9461
$ws.Range("M134").Value = -3169.9461

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 810.34784
$ws.Range("I7").Value = 1002.1177
$ws.Range("K7").Value = 3006.3531
$ws.Range("M7").Value = -2894.3531
$ws.Range("H92").Value = 999
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2997
$ws.Range("N92").Value = -5493
$ws.Range("M92").ClearContents()
$ws.Range("H107").Value = 650
$ws.Range("J107").Value = 650
$ws.Range("L107").Value = 1950
$ws.Range("N107").Value = -5790

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 412.3889
$ws.Range("I97").Value = 422.2
$ws.Range("J97").Value = 363.33334
$ws.Range("K97").Value = 422.2
$ws.Range("L97").Value = 363.33334
$ws.Range("M97").Value = 73.80000000000001
$ws.Range("N97").Value = -1355.33334
$ws.Range("H132").Value = 2423.4546
$ws.Range("I132").Value = 2060.8
$ws.Range("K132").Value = 6182.400000000001
$ws.Range("M132").Value = -3652.400000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5201
$ws.Range("I7").Value = 2785.6667
$ws.Range("J7").Value = 6006.1113
$ws.Range("K7").Value = 2785.6667
$ws.Range("L7").Value = 6006.1113
$ws.Range("M7").Value = -2673.6667
$ws.Range("N7").Value = -6230.1113
$ws.Range("H55").Value = 1198.238
$ws.Range("I55").Value = 1264.6666
$ws.Range("K55").Value = 1264.6666
$ws.Range("M55").Value = -1091.6666
$ws.Range("H82").Value = 4963.8423
$ws.Range("I82").Value = 2683
$ws.Range("J82").Value = 8100
$ws.Range("K82").Value = 2683
$ws.Range("L82").Value = 8100
$ws.Range("M82").Value = -2322
$ws.Range("N82").Value = -8822
$ws.Range("H85").Value = 4963.8423
$ws.Range("I85").Value = 2683
$ws.Range("J85").Value = 8100
$ws.Range("K85").Value = 2683
$ws.Range("L85").Value = 8100
$ws.Range("M85").Value = -1435
$ws.Range("N85").Value = -10596
$ws.Range("H93").Value = 6467.75
$ws.Range("I93").Value = 6467.75
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 6467.75
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -5219.75
$ws.Range("N93").ClearContents()
$ws.Range("H118").Value = 41213.93
$ws.Range("J118").Value = 41213.93
$ws.Range("L118").Value = 41213.93
$ws.Range("N118").Value = -44527.93
$ws.Range("H126").Value = 5201
$ws.Range("I126").Value = 2785.6667
$ws.Range("J126").Value = 6006.1113
$ws.Range("K126").Value = 8357.000100000001
$ws.Range("L126").Value = 18018.3339
$ws.Range("M126").Value = -5887.000100000001
$ws.Range("N126").Value = -22958.3339

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 149
$ws.Range("I107").Value = 149
$ws.Range("K107").Value = 447
$ws.Range("M107").Value = 1473
$ws.Range("H122").Value = 2743.8
$ws.Range("I122").Value = 2307.9119
$ws.Range("K122").Value = 6923.7357
$ws.Range("M122").Value = -4473.7357
$ws.Range("H136").Value = 1780.0197
$ws.Range("J136").Value = 3486.4443
$ws.Range("L136").Value = 10459.3329
$ws.Range("N136").Value = -15559.3329
